# Append: 2025-12-12 01:56 JST
# The scraper re-ran and produced a new batch with the same rows but a
# fresh acquisition timestamp. Update the "取得日時" (acquisition
# timestamp) column (A) for all existing data rows on the active sheet
# ("ランサーズ") from the previous run's timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2025-12-12 01:24:35"
$newTimestamp = "2025-12-12 01:56:22"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
